$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.333947
$ws.Range("H2").Value = 1.001841
$ws.Range("I2").Value = 0.9184595666969813
$ws.Range("J2").Value = 0.9184595666969813
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07259900000000001
$ws.Range("N2").Value = 0.217797
$ws.Range("O2").Value = 0.0162094769588191
$ws.Range("P2").Value = 0.0162094769588191
$ws.Range("Q2").Value = 0.024244218253
$ws.Range("R2").Value = 0.218197964277
$ws.Range("S2").Value = 0.01488774918398169
$ws.Range("T2").Value = 0.01488774918398169
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.333947
$ws.Range("H3").Value = 1.001841
$ws.Range("I3").Value = 0.9184595666969813
$ws.Range("J3").Value = 0.9184595666969813
$ws.Range("O3").Value = 0.9349228167457665
$ws.Range("P3").Value = 0.9349228167457664
$ws.Range("Q3").Value = 1.398346959404
$ws.Range("R3").Value = 12.585122634636
$ws.Range("S3").Value = 0.8586888051634379
$ws.Range("T3").Value = 0.8586888051634378
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.333947
$ws.Range("H4").Value = 1.001841
$ws.Range("I4").Value = 0.9184595666969813
$ws.Range("J4").Value = 0.9184595666969813
$ws.Range("M4").Value = 0.2188686666666667
$ws.Range("N4").Value = 0.656606
$ws.Range("O4").Value = 0.04886770629541442
$ws.Range("P4").Value = 0.04886770629541441
$ws.Range("Q4").Value = 0.07309053462733334
$ws.Range("R4").Value = 0.657814811646
$ws.Range("S4").Value = 0.04488301234956167
$ws.Range("T4").Value = 0.04488301234956166
$ws.Range("I5").Value = 0.08154043330301874
$ws.Range("J5").Value = 0.08154043330301874
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.07259900000000001
$ws.Range("N5").Value = 0.217797
$ws.Range("O5").Value = 0.0162094769588191
$ws.Range("P5").Value = 0.0162094769588191
$ws.Range("Q5").Value = 0.002152390952333333
$ws.Range("R5").Value = 0.019371518571
$ws.Range("S5").Value = 0.001321727774837408
$ws.Range("T5").Value = 0.001321727774837408
$ws.Range("I6").Value = 0.08154043330301874
$ws.Range("J6").Value = 0.08154043330301874
$ws.Range("O6").Value = 0.9349228167457665
$ws.Range("P6").Value = 0.9349228167457664
$ws.Range("S6").Value = 0.07623401158232859
$ws.Range("T6").Value = 0.07623401158232858
$ws.Range("I7").Value = 0.08154043330301874
$ws.Range("J7").Value = 0.08154043330301874
$ws.Range("M7").Value = 0.2188686666666667
$ws.Range("N7").Value = 0.656606
$ws.Range("O7").Value = 0.04886770629541442
$ws.Range("P7").Value = 0.04886770629541441
$ws.Range("Q7").Value = 0.006488945273111111
$ws.Range("R7").Value = 0.058400507458
$ws.Range("S7").Value = 0.003984693945852749
$ws.Range("T7").Value = 0.003984693945852748
